$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Hora: 17:00 às 18:00" -> "Hora: 21:00 às 22:30"
#    The new text is split across several runs with varying rPr, so we
#    rebuild the whole paragraph via InsertXML using the exact OOXML
#    that the target revision contains (preserving the paragraph's
#    pPr and each run's rPr).
# ------------------------------------------------------------------
$t = $d.Tables.Item(1)
$horaRow = 0
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $cellText = $t.Cell($i, 1).Range.Text
    if ($cellText -like "Hora:*") {
        $horaRow = $i
        break
    }
}

if ($horaRow -gt 0) {
    $cell = $t.Cell($horaRow, 1)
    $para = $cell.Range.Paragraphs.Item(1)
    $r = $para.Range

    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:widowControl w:val="false"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:ascii="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Hora: </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="Times New Roman" w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="pt-BR" w:eastAsia="pt-BR" w:bidi="ar-SA"/></w:rPr><w:t>21:</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:ascii="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">00 às </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:ascii="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>22</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:ascii="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:ascii="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:ascii="Calibri" w:hAnsi="Calibri"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>0</w:t></w:r></w:p>'

    $r.InsertXML($xml)
}

# ------------------------------------------------------------------
# 2) "Google Meet" -> "Teams" (simple text replace, same formatting)
# ------------------------------------------------------------------
$d.Content.Find.Execute("Google Meet", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Teams", 2)
